# KPI update: refresh "last_edited_time" (column D) timestamps
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lũy kế ngày LONG XUYÊN")

# Rows 2-26 -> 2024-07-25T16:53:00.000Z
$ws.Range("D2:D26").Value = "2024-07-25T16:53:00.000Z"

# Rows 27-84 -> 2024-07-25T16:54:00.000Z
$ws.Range("D27:D84").Value = "2024-07-25T16:54:00.000Z"

# Rows 85-88 -> 2024-07-25T16:55:00.000Z
$ws.Range("D85:D88").Value = "2024-07-25T16:55:00.000Z"
